# Manage news test cases
# Add a new LoginPage test-data row ("admin1" / "admin") and leave the
# LoginPage sheet as the active/selected sheet with E7 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginPage")

# Bring LoginPage to the front (it becomes the workbook's tabSelected /
# active sheet instead of HomePage).
$ws.Activate()

# New row of login test data: Username="admin1", password="admin".
$ws.Cells.Item(7, 1).Value = "admin1"
$ws.Cells.Item(7, 2).Value = "admin"

# Leave the selection where the author last clicked.
$ws.Range("E7").Select()
